$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update existing cell values in column C and D (rows 23-75) ---
$ws.Range("D23").Value = 0
$ws.Range("C26").Value = 80.90000000000001
$ws.Range("D30").Value = 0.4
$ws.Range("D32").Value = 2.3
$ws.Range("D51").Value = 1.4
$ws.Range("C55").Value = 105.5
$ws.Range("D55").Value = -0.1
$ws.Range("C58").Value = 105.1
$ws.Range("D59").Value = 1.2
$ws.Range("C60").Value = 108.2
$ws.Range("C61").Value = 109
$ws.Range("D61").Value = 0.7
$ws.Range("C62").Value = 110.3
$ws.Range("C63").Value = 111.6
$ws.Range("D63").Value = 1.2
$ws.Range("C65").Value = 112
$ws.Range("D65").Value = 0.9
$ws.Range("C66").Value = 111.8
$ws.Range("D66").Value = -0.1
$ws.Range("C67").Value = 113.6
$ws.Range("D67").Value = 1.5
$ws.Range("C68").Value = 114.1
$ws.Range("D68").Value = 0.4
$ws.Range("C69").Value = 109.3
$ws.Range("D69").Value = -4.2
$ws.Range("C70").Value = 111.7
$ws.Range("D70").Value = 2.2
$ws.Range("C71").Value = 98.09999999999999
$ws.Range("D71").Value = -12.2
$ws.Range("C72").Value = 102.9
$ws.Range("D72").Value = 4.9
$ws.Range("C73").Value = 109.1
$ws.Range("D73").Value = 6
$ws.Range("C74").Value = 113.2
$ws.Range("D74").Value = 3.7
$ws.Range("C75").Value = 115.1
$ws.Range("D75").Value = 1.7

# --- Append new row 76 with the latest quarter data ---
# Force the date-like label to be stored as text (matches the other
# "Serie" column cells, which are shared-string text, not dates).
$ws.Range("A76").NumberFormat = "@"
$ws.Range("A76").Value = "01-07-2021"
$ws.Range("A76").Style = "Normal"

$ws.Range("B76").Value = 117.7
$ws.Range("C76").Value = 121.2
$ws.Range("D76").Value = 5.3
